# Apply edits described by the diff to empresa_banco_dados.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet: Empregado ---
$ws = $wb.Worksheets.Item("Empregado")

$ws.Range("A1").Value = "nss (chave)"

$ws.Range("A2").Value = 1001
$ws.Range("B2").Value = "João"
$ws.Range("C2").Value = "Carlos"
$ws.Range("D2").Value = "Silva"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1990-05-15"
$ws.Range("G2").Value = 3500
$ws.Range("H2").Value = "Rua A, 123"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5001

$ws.Range("A3").Value = 1002
$ws.Range("B3").Value = "Maria"
$ws.Range("C3").Value = "José"
$ws.Range("D3").Value = "Souza"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "1992-08-22"
$ws.Range("G3").Value = 3200
$ws.Range("H3").Value = "Rua B, 456"
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 5002

# --- Sheet: Departamento ---
$ws = $wb.Worksheets.Item("Departamento")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "TI"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 1001
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2020-01-10"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "RH"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1002
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2021-02-15"

# --- Sheet: Projeto ---
$ws = $wb.Worksheets.Item("Projeto")

$ws.Range("A2").Value = 101
$ws.Range("B2").Value = "Projeto A"
$ws.Range("D2").Value = 1

$ws.Range("A3").Value = 102
$ws.Range("B3").Value = "Projeto B"
$ws.Range("D3").Value = 2

# --- Sheet: Depende ---
$ws = $wb.Worksheets.Item("Depende")

$ws.Range("A2").Value = 1001
$ws.Range("B2").Value = "Maria Silva"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2015-07-20"
$ws.Range("E2").Value = "Filha"

$ws.Range("A3").Value = 1002
$ws.Range("B3").Value = "João Souza"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2018-11-10"
$ws.Range("E3").Value = "Filho"

# --- Sheet: Trabalha-em ---
$ws = $wb.Worksheets.Item("Trabalha-em")

$ws.Range("A2").Value = 1001
$ws.Range("B2").Value = 101
$ws.Range("C2").Value = 40

$ws.Range("A3").Value = 1002
$ws.Range("B3").Value = 102
$ws.Range("C3").Value = 35

# --- Sheet: Localizacao ---
$ws = $wb.Worksheets.Item("Localizacao")

$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
